# Fruta / hortaliza, semanal
# Update the weekly date/volume/price figures for rows 2-10 (Guayaba, Vega
# Modelo de Temuco) to the new values for this week's refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = fecha; M = volumen; N = precio minimo; O = precio maximo; P = precio promedio ponderado; S = precio $/kg }
$updates = @{
    2  = @{ D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    3  = @{ D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    4  = @{ D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    5  = @{ D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
    6  = @{ D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 }
    7  = @{ D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 }
    9  = @{ D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
    10 = @{ D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
